# Updated cryptos list on Fri May 10 13:15:03 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.967.26"
$ws.Range("E2").Value = "  +2.43%  "
$ws.Range("D3").Value = "3.028.86"
$ws.Range("E3").Value = "  +1.40%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  +0.08%  "
$__style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.00"
$ws.Range("D6").Style = $__style
$ws.Range("E6").Value = "  +6.03%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.024.43"
$ws.Range("E8").Value = "  +1.30%  "
$__style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.513"
$ws.Range("D9").Style = $__style
$ws.Range("E9").Value = "  -0.29%  "
$__style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.11"
$ws.Range("D10").Style = $__style
$ws.Range("E10").Value = "  +17.23%  "
$ws.Range("E11").Value = "  +1.81%  "
$ws.Range("E12").Value = "  +2.15%  "
$ws.Range("E13").Value = "  +2.87%  "
$__style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.60"
$ws.Range("D14").Style = $__style
$ws.Range("E14").Value = "  +4.16%  "
$__style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.124"
$ws.Range("D15").Style = $__style
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("D16").Value = "3.531.26"
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$__style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.07"
$ws.Range("D17").Style = $__style
$ws.Range("E17").Value = "  +2.63%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "62.967.76"
$ws.Range("E18").Value = "  +2.51%  "
$ws.Range("D19").Value = "3.030.51"
$ws.Range("E19").Value = "  +1.49%  "
$__style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "448.90"
$ws.Range("D20").Style = $__style
$ws.Range("E20").Value = "  -0.10%  "
$__style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.23"
$ws.Range("D21").Style = $__style
$ws.Range("E21").Value = "  +1.56%  "
$ws.Range("E22").Value = "  +1.95%  "
$ws.Range("E23").Value = "  +2.74%  "
$__style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.48"
$ws.Range("D24").Style = $__style
$ws.Range("E24").Value = "  +8.59%  "
$__style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.16"
$ws.Range("D25").Style = $__style
$ws.Range("E25").Value = "  +1.63%  "
$__style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.32"
$ws.Range("D26").Style = $__style
$ws.Range("E26").Value = "  +6.58%  "
$ws.Range("E27").Value = "  +3.12%  "
$ws.Range("E28").Value = "  +0.02%  "
$__style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.57"
$ws.Range("D29").Style = $__style
$ws.Range("E29").Value = "  +5.04%  "
$__style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.27"
$ws.Range("D30").Style = $__style
$ws.Range("E30").Value = "  +10.69%  "
$__style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.70"
$ws.Range("D31").Style = $__style
$ws.Range("E31").Value = "  +0.76%  "
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("E33").Value = "  +1.82%  "
$ws.Range("E34").Value = "  +1.60%  "
$ws.Range("D35").Value = "0.0₃0878"
$ws.Range("E35").Value = "  +6.82%  "
$__style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.05"
$ws.Range("D36").Style = $__style
$ws.Range("E36").Value = "  +3.36%  "
$__style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.88"
$ws.Range("D37").Style = $__style
$ws.Range("E37").Value = "  +1.71%  "
$__style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.17"
$ws.Range("D38").Style = $__style
$ws.Range("E38").Value = "  +10.26%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$__style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.11"
$ws.Range("D39").Style = $__style
$ws.Range("E39").Value = "  +2.97%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$__style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.130"
$ws.Range("D40").Style = $__style
$ws.Range("E40").Value = "  +8.42%  "
$__style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.51"
$ws.Range("D41").Style = $__style
$ws.Range("E41").Value = "  +0.59%  "
$__style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.06"
$ws.Range("D42").Style = $__style
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("B43").Value = "Arweave"
$ws.Range("C43").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$__style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "44.53"
$ws.Range("D43").Style = $__style
$ws.Range("E43").Value = "  +16.60%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$__style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.306"
$ws.Range("D44").Style = $__style
$ws.Range("E44").Value = "  +13.95%  "
$__style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "393.11"
$ws.Range("D45").Style = $__style
$ws.Range("E45").Value = "  +1.23%  "
$__style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0359"
$ws.Range("D46").Style = $__style
$ws.Range("E46").Value = "  +2.42%  "
$ws.Range("D47").Value = "2.712.42"
$ws.Range("E47").Value = "  +0.67%  "
$__style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.77"
$ws.Range("D48").Style = $__style
$ws.Range("E48").Value = "  +2.15%  "
$__style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.49"
$ws.Range("D49").Style = $__style
$ws.Range("E49").Value = "  +13.36%  "
$__style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.27"
$ws.Range("D51").Style = $__style
$ws.Range("E51").Value = "  +5.88%  "
